# Update the Initial Weights, Opt Portfolio, and Opt Portfolio with View
# columns to reflect the recalculated equal-weighted initial allocation
# (1/7 per asset) and the resulting re-optimized portfolio weights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(0.1428571428571428, 0.1521569877763614, 0.1521569873121236)
    3 = @(0.1428571428571428, 0.1514628155394849, 0.1514628154088595)
    4 = @(0.1428571428571428, 0.1307254289492208, 0.1307254287537006)
    5 = @(0.1428571428571428, 0.1307357372099624, 0.130735737147938)
    6 = @(0.1428571428571428, 0.1344673604747327, 0.1344673608767446)
    7 = @(0.1428571428571428, 0.1501043822555428, 0.1501043828279517)
    8 = @(0.1428571428571428, 0.1503472877946951, 0.150347287672682)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]
    $ws.Cells.Item($row, 3).Value = $rowVals[1]
    $ws.Cells.Item($row, 4).Value = $rowVals[2]
}
